# Natmi following Dr Hou advice
# Expand the single FAPs->M2 (Penk-Oprm1) edge into the full set of
# sending/target cluster combinations (FAPs, M1, sCs) x (M1, M2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("FAPs","Penk","Oprm1","M1",3,1,20.908407,62.725221,0.9843255515709016,0.9867193626172466,1,0.3333333333333333,0.05215766666666666,0.156473,0.1010355835763341,0.1010355835763341,1.090533722837,9.814803505533,0.09945190653206297,0.0996937666281019),
    @("FAPs","Penk","Oprm1","M2",3,1,20.908407,62.725221,0.9843255515709016,0.9867193626172466,3,1,0.464073,1.392219,0.8989644164236659,0.8989644164236659,9.703027161711,87.327244455399,0.8848736450388386,0.8870255959891447),
    @("M1","Penk","Oprm1","M1",1,0.3333333333333333,0.17835,0.53505,0.00839635760498972,0.008416776960711829,1,0.3333333333333333,0.05215766666666666,0.156473,0.1010355835763341,0.1010355835763341,0.00930231985,0.08372087865000001,0.0008483308905357271,0.0008503939720573631),
    @("M1","Penk","Oprm1","M2",1,0.3333333333333333,0.17835,0.53505,0.00839635760498972,0.008416776960711829,3,1,0.464073,1.392219,0.8989644164236659,0.8989644164236659,0.08276741955000001,0.7449067759500001,0.007548026714453993,0.007566382988654465),
    @("sCs","Penk","Oprm1","M1",1,0.5,0.1545965,0.309193,0.007278090824108737,0.004863860422041627,1,0.3333333333333333,0.05215766666666666,0.156473,0.1010355835763341,0.1010355835763341,0.008063392714833332,0.048380356289,0.0007353461537353885,0.0004914229761748103),
    @("sCs","Penk","Oprm1","M2",1,0.5,0.1545965,0.309193,0.007278090824108737,0.004863860422041627,3,1,0.464073,1.392219,0.8989644164236659,0.8989644164236659,0.0717440615445,0.430464369267,0.006542744670373348,0.004372437445866816)
)

$r = 2
foreach ($rowdata in $rows) {
    $c = 1
    foreach ($val in $rowdata) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
